# Generate Report for Handoff
#
# The handoff/handback polling run that produced this report re-ran after
# the previous snapshot, advancing the "Handback transform failed" row
# (4ca241b1-...) and the "Ready for handoff" rows (2e874e4a-..., 66d791fb-...,
# 6cff4441-..., e425c7da-..., f56d6f71-..., fae5ebb3-...) to a newer
# handoff timestamp. Update the recorded datetimes on every sheet to match
# the latest run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("D7").Value  = "2016-23-20 06:23:31"
$ws.Range("D10").Value = "2016-23-20 06:23:31"
$ws.Range("D11").Value = "2016-23-20 06:23:31"
$ws.Range("D12").Value = "2016-23-20 06:23:31"
$ws.Range("D13").Value = "2016-23-20 06:23:31"
$ws.Range("D14").Value = "2016-23-20 06:23:31"
$ws.Range("D15").Value = "2016-23-20 06:23:31"
$ws.Range("D16").Value = "2016-23-20 06:23:31"

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("E7").Value  = "2016-03-20 06:23:27"
$ws.Range("E10").Value = "2016-03-20 06:23:27"
$ws.Range("E11").Value = "2016-03-20 06:23:27"
$ws.Range("E12").Value = "2016-03-20 06:23:27"
$ws.Range("E13").Value = "2016-03-20 06:23:27"
$ws.Range("E14").Value = "2016-03-20 06:23:27"
$ws.Range("E15").Value = "2016-03-20 06:23:27"
$ws.Range("E16").Value = "2016-03-20 06:23:27"

# --- de-de sheet: "Latest Handoff Datetime" column (E) ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("E7").Value  = "2016-03-20 06:23:31"
$ws.Range("E10").Value = "2016-03-20 06:23:31"
$ws.Range("E11").Value = "2016-03-20 06:23:31"
$ws.Range("E12").Value = "2016-03-20 06:23:31"
$ws.Range("E13").Value = "2016-03-20 06:23:31"
$ws.Range("E14").Value = "2016-03-20 06:23:31"
$ws.Range("E15").Value = "2016-03-20 06:23:31"
$ws.Range("E16").Value = "2016-03-20 06:23:31"
